# Adapt column header formatting to respective input file names (#7)
#
# Renames the shared "_old"/"_new" header-name suffixes to the concrete
# format-version identifiers ("_FV2304" / "_FV2310"), turns the sheet's
# data range into a native Excel Table, and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename header cells (row 1) from "<Name>_old" / "<Name>_new" to
#    "<Name>_FV2304" / "<Name>_FV2310" respectively.
# ---------------------------------------------------------------------
$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count
$headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $lastCol))
for ($i = 1; $i -le $headerRange.Columns.Count; $i++) {
    $cell = $headerRange.Cells.Item(1, $i)
    $value = $cell.Value2
    if ($value -ne $null) {
        if ($value.EndsWith("_old")) {
            $cell.Value = $value.Substring(0, $value.Length - 4) + "_FV2304"
        } elseif ($value.EndsWith("_new")) {
            $cell.Value = $value.Substring(0, $value.Length - 4) + "_FV2310"
        }
    }
}

# ---------------------------------------------------------------------
# 2) Convert the used range into an Excel Table (ListObject) so the
#    renamed headers double as table column headers.
# ---------------------------------------------------------------------
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $usedRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# ---------------------------------------------------------------------
# 3) Freeze the header row (split/freeze pane under row 1).
# ---------------------------------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$null
